# Center and underline the document's title heading
# ("GDI (Graphics Device Interface) - Summary for VCL"), which is the
# very first paragraph of the document.
$d = $word.ActiveDocument

$title = $d.Paragraphs(1)

# Center the paragraph (adds <w:jc w:val="center"/> to the pPr).
$title.Alignment = 1

# Underline the whole paragraph, including its paragraph mark, so both
# the paragraph-mark run properties and the text run properties get
# <w:u w:val="single"/>.
$title.Range.Font.Underline = 1
